# Update column G ("K") values on the active worksheet to reflect the
# regenerated save_data (K instead of Strike#, regen std/mean, calc and
# write s_vals).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 1
    3  = 1
    4  = 1
    5  = 0
    6  = 1
    7  = 4
    8  = 3
    9  = 2
    10 = 0
    11 = 0
    13 = 0
    14 = 3
    15 = 3
    16 = 2
    17 = 1
    18 = 0
    19 = 1
    20 = 2
    21 = 2
    22 = 2
    23 = 0
    25 = 0
    26 = 1
    27 = 1
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
